$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new "Save" column
$ws.Cells.Item(1, 8).Value = "Save"
$ws.Range("G1").Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)

# Save values per row (row number => value)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
